# #27 finish bootcamp join in dimension
# Fill in the "joinDate" column (H) on the "privateinfo" sheet with the
# bootcamp join-window timestamps used by the test fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("privateinfo")
$ws.Activate()

$ws.Range("H2").Value = "2022-5-1 00:00:00"
$ws.Range("H3").Value = "2022-5-1 00:00:00"
$ws.Range("H4").Value = "2022-5-1 00:00:00"
$ws.Range("H5").Value = "2022-5-1 00:00:00"
$ws.Range("H6").Value = "2022-5-1 00:00:00"
$ws.Range("H7").Value = "2022-5-1 23:59:59"
$ws.Range("H8").Value = "2022-5-1 23:59:59"
$ws.Range("H9").Value = "2022-5-1 23:59:59"
$ws.Range("H10").Value = "2022-5-2 00:00:00"
$ws.Range("H11").Value = "2022-5-3 00:00:01"
$ws.Range("H12").Value = "2022-5-4 00:00:02"
$ws.Range("H13").Value = "2022-5-5 00:00:03"
$ws.Range("H14").Value = "2022-5-6 00:00:04"
$ws.Range("H15").Value = "2022-5-7 00:00:05"
$ws.Range("H16").Value = "2022-5-1 23:59:59"
$ws.Range("H17").Value = "2022-5-2 00:00:00"
$ws.Range("H18").Value = "2022-5-3 00:00:01"
$ws.Range("H19").Value = "2022-5-4 00:00:02"
$ws.Range("H20").Value = "2022-5-5 00:00:03"
$ws.Range("H21").Value = "2022-5-6 00:00:04"
$ws.Range("H22").Value = "2022-5-7 00:00:05"
$ws.Range("H23").Value = "2022-5-1 23:59:59"
$ws.Range("H24").Value = "2022-5-2 00:00:00"
$ws.Range("H25").Value = "2022-5-3 00:00:01"
$ws.Range("H26").Value = "2022-5-4 00:00:02"
$ws.Range("H27").Value = "2022-5-5 00:00:03"
$ws.Range("H28").Value = "2022-5-6 00:00:04"
$ws.Range("H29").Value = "2022-5-7 00:00:05"
$ws.Range("H30").Value = "2022-5-1 23:59:59"
$ws.Range("H31").Value = "2022-5-2 00:00:00"
$ws.Range("H32").Value = "2022-5-3 00:00:01"
$ws.Range("H33").Value = "2022-5-4 00:00:02"
$ws.Range("H34").Value = "2022-5-5 00:00:03"
$ws.Range("H35").Value = "2022-5-6 00:00:04"
$ws.Range("H36").Value = "2022-5-7 00:00:05"
$ws.Range("H37").Value = "2022-5-1 23:59:59"
$ws.Range("H38").Value = "2022-5-2 00:00:00"
$ws.Range("H39").Value = "2022-5-3 00:00:01"
$ws.Range("H40").Value = "2022-5-4 00:00:02"
$ws.Range("H41").Value = "2022-5-5 00:00:03"
$ws.Range("H42").Value = "2022-5-6 00:00:04"
$ws.Range("H43").Value = "2022-5-7 00:00:05"
$ws.Range("H44").Value = "2022-5-6 00:00:04"

# Mirror the author's final selection/scroll state: cursor left on H41,
# with the view no longer pinned to a manually scrolled top-left cell.
[void]$ws.Range("A1").Select()
[void]$ws.Range("H41").Select()
